$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("To do")
$cfs = $ws.Cells.FormatConditions
Write-Output ("before count: " + $cfs.Count())
for ($i=1; $i -le $cfs.Count(); $i++) {
    $fc = $cfs.Item($i)
    Write-Output ($i.ToString() + " | " + $fc.AppliesTo().Address() + " | " + $fc.Formula1())
}
$ws.Rows.Item(20).Insert()
Write-Output "---after insert---"
$cfs2 = $ws.Cells.FormatConditions
Write-Output ("after count: " + $cfs2.Count())
for ($i=1; $i -le $cfs2.Count(); $i++) {
    $fc = $cfs2.Item($i)
    Write-Output ($i.ToString() + " | " + $fc.AppliesTo().Address() + " | " + $fc.Formula1())
}
